# Appends one new candidate (陈照明) to the Candidates sheet as row 12,
# mirroring the layout/formatting of the most recently-added candidate
# row (row 11), wiring up the e-mail hyperlink and restoring the wrap
# height on row 9 that Excel recalculates once the new row is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12: clone formatting from row 11 (the previous "last" row) ---
$ws.Range("A11:W11").Copy()
$ws.Range("A12:W12").PasteSpecial(-4122)  # xlPasteFormats

# --- Candidate data ---
$ws.Range("A12").Value = "陈照明 "
$ws.Range("B12").Value = "zhaoming3117@qq.com"
$ws.Range("C12").Value = "中"
$ws.Range("D12").Value = 43583
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "大专"
$ws.Range("G12").Value = 33
$ws.Range("H12").Value = "10+"
$ws.Range("I12").Value = "室内设计"
$ws.Range("J12").Value = "已收到"
$ws.Range("K12").Value = "BOSS直聘"
$ws.Range("L12").Value = "N"
$ws.Range("N12").Value = "N"
$ws.Range("P12").Value = "N"
$ws.Range("Q12").Value = "N"
$ws.Range("R12").Value = "N"

# E-mail hyperlink for the new candidate, same convention as B9/B10/B11.
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:zhaoming3117@qq.com")

# Hyperlinks.Add() mints a fresh (near-duplicate) cell style; re-stamp B12
# with B11's existing hyperlink style so it matches the rest of the column.
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)  # xlPasteFormats

# Row heights: the new row settles at 16.5pt; row 9's wrapped description
# reflows to 49.5pt once the sheet is resaved.
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 49.5

# Keep the active selection near the newly appended row.
$ws.Range("B30").Select()
